$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# The "Name" row (A4) loses its value text, which moves to become the
# "Title" row's (A5) value instead of the old placeholder title text.
$ws.Range("B4").Value = ""
$ws.Range("B5").Value = 'Mapping Métier/CDA/FHIR : "Directive Anticipee"'

# Bump the recorded generation date/time.
$ws.Range("B8").Value = "2026-01-07T15:20:53+00:00"
